# market_health_data.xlsx update: 2025-10-29 18:32
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Metadata sheet: bump "Last Updated" timestamp by one minute.
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "29 Oct 2025, 06:32 PM"

# ---------------------------------------------------------------------
# 2) Top Gainers sheet: refreshed ranking for rows 63-66.
#    NPST jumps to the top of this block; ORIENTTECH/ICRA/SALASAR each
#    shift down one spot, re-using their existing (unchanged) figures.
# ---------------------------------------------------------------------
$wsGainers = $wb.Worksheets.Item("Top Gainers")

$gainersRows = @(
    @(63, "NPST",       3.8509, -2.0059, -3.5057),
    @(64, "ORIENTTECH", 3.827,   0.5247000000000001, 32.6784),
    @(65, "ICRA",        3.7985, 4.4793, 2.8828),
    @(66, "SALASAR",     3.7935, 4.7872, 11.0485)
)

foreach ($row in $gainersRows) {
    $r = $row[0]
    $wsGainers.Cells.Item($r, 2).Value = $row[1]
    $wsGainers.Cells.Item($r, 3).Value = $row[2]
    $wsGainers.Cells.Item($r, 4).Value = $row[3]
    $wsGainers.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------
# 3) Top Losers sheet: refreshed ranking for rows 51-76.
#    UNIMECH and ALLDIGI move up with newly computed figures, JNKINDIA
#    is newly added at row 71, everything else shifts down one spot,
#    and BBOX (previously row 76) drops off the bottom of the list.
# ---------------------------------------------------------------------
$wsLosers = $wb.Worksheets.Item("Top Losers")

$losersRows = @(
    @(51, "UNIMECH",    -2.8008, -1.6104, -0.4585),
    @(52, "TTKPRESTIG", -2.7438, 8.001200000000001, 9.650499999999999),
    @(53, "PFOCUS",     -2.7039, -2.6276, -1.2163),
    @(54, "ALLDIGI",    -2.6342, -0.2306, -5.3103),
    @(55, "PRIVISCL",   -2.6288, -2.1048, 19.7451),
    @(56, "CANHLIFE",   -2.6148, "N/A", "N/A"),
    @(57, "GKENERGY",   -2.6122, -9.807700000000001, 23.2758),
    @(58, "SGFIN",      -2.592,  -0.06270000000000001, 11.7235),
    @(59, "ARVINDFASN", -2.549,  -2.9892, -4.4223),
    @(60, "EDELWEISS",  -2.5422, -3.3745, 8.5305),
    @(61, "SAMHI",      -2.5284, 1.8231, 2.8516),
    @(62, "TBOTEK",     -2.524,  -3.5732, 1.036),
    @(63, "UJJIVANSFB", -2.5201, 0.3845, 12.6645),
    @(64, "AMBER",      -2.5098, -0.1082, 2.763),
    @(65, "GRPLTD",     -2.4898, -5.9894, -5.4586),
    @(66, "NESCO",      -2.4722, 1.9934, 3.8931),
    @(67, "PILANIINVS", -2.4546, -0.7907, 4.267),
    @(68, "NSIL",       -2.4088, -1.7646, 4.7431),
    @(69, "COALINDIA",  -2.4016, -3.058, -2.0387),
    @(70, "FINOPB",     -2.3673, -6.2696, 11.1938),
    @(71, "JNKINDIA",   -2.3482, -2.8371, 4.2622),
    @(72, "FCL",        -2.3453, -2.616, -0.02),
    @(73, "DEEDEV",     -2.3334, -6.6528, -7.4227),
    @(74, "WEALTH",     -2.3047, -3.8606, -2.8234),
    @(75, "RATNAMANI",  -2.2788, -0.4626, 0.8712),
    @(76, "CSBBANK",    -2.2695, 2.3137, 10.6999)
)

foreach ($row in $losersRows) {
    $r = $row[0]
    $wsLosers.Cells.Item($r, 2).Value = $row[1]
    $wsLosers.Cells.Item($r, 3).Value = $row[2]
    $wsLosers.Cells.Item($r, 4).Value = $row[3]
    $wsLosers.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------
# 4) "distance from Dma50" sheet: the three trailing N/A sector rows
#    (NIFTYFINSEREXBNK, NIFTYMSITTELCM, NIFTYMSFINSERV) are removed.
# ---------------------------------------------------------------------
$wsDma50 = $wb.Worksheets.Item("distance from Dma50")
$wsDma50.Rows("31:33").Delete()
